# Auto-generated edit script applying the cryptos.xlsx price/volume update
# Commit: Updated cryptos list on Mon Sep 25 06:38:53 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.215.07'
$ws.Range("E2").Value = '  -1.90%  '
$ws.Range("D3").Value = '1.581.91'
$ws.Range("E3").Value = '  -1.22%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = '209.51'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("E6").Value = '  -2.48%  '
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").Value = '19.52'
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '1.802.93'
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("D13").Value = '1.573.47'
$ws.Range("E13").Value = '  -1.77%  '
$ws.Range("D14").Value = '4.03'
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").Value = '0.517'
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("D16").Value = '64.40'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").Value = '26.206.57'
$ws.Range("E17").Value = '  -1.82%  '
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").Value = '7.25'
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").Value = '206.23'
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("D22").Value = '4.26'
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").Value = '2.21'
$ws.Range("E23").Value = '  -2.61%  '
$ws.Range("D24").Value = '8.85'
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("D25").Value = '144.83'
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").Value = '7.01'
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("D29").Value = '15.23'
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("D33").Value = '2.94'
$ws.Range("E33").Value = '  -1.12%  '
$ws.Range("D34").Value = '1.283.23'
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("E35").Value = '  +7.98%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("E39").Value = '  -1.92%  '
$ws.Range("D40").Value = '0.813'
$ws.Range("E40").Value = '  -1.85%  '
$ws.Range("D41").Value = '5.55'
$ws.Range("E41").Value = '  +2.61%  '
$ws.Range("D42").Value = '0.769'
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("D44").Value = '62.12'
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("D45").Value = '1.716.17'
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("D46").Value = '88.57'
$ws.Range("E46").Value = '  -2.39%  '
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  -0.53%  '
$ws.Range("E49").Value = '  -1.51%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₇0955'
$ws.Range("E51").Value = '  -10.16%  '
